$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the developer vacation date intervals: append an extra interval for
# Alice and add trailing "|" separators for Charlie and Dave. Charlie/Dave are
# updated before Alice so the shared-string table keeps the same relative
# ordering that Excel itself produces.
$ws.Range("C4").Value = "2025-04-05;2025-04-10|"
$ws.Range("C5").Value = "2025-02-25;2025-02-28|"
$ws.Range("C2").Value = "2025-03-10;2025-03-15|2025-05-10;2025-05-15"

# Turn the data range into a proper Excel Table ("Table1").
$listObject = $ws.ListObjects.Add(1, $ws.Range("A1:D5"), $null, 1)
$listObject.Name = "Table1"

# Resize columns to fit the new table look (drop bestFit autosizing).
$ws.Columns.Item(1).ColumnWidth = 12.385416666666666
$ws.Columns.Item(2).ColumnWidth = 16.721354166666668
$ws.Columns.Item(4).ColumnWidth = 16.944010416666668

# Update the active selection to match the saved view state.
$null = $ws.Range("C3").Select()
